# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-10
$newValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 2
    7  = 1
    8  = 1
    9  = 1
    10 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
